$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Content edits (do these before sheet deletion / selection changes so
#    shared-string bookkeeping settles naturally). The order of first-use
#    below matches the order new shared strings appear in the saved file.
# ---------------------------------------------------------------------------

# deleteUser sheet: rename the "deleteAdminNotPossibleUser" TC.
$wsDeleteUser = $wb.Worksheets.Item("deleteUser")
$wsDeleteUser.Range("A2").Value = "deleteAdminNotPossible"

# modifyUser sheet: the "modifyUser_Passwd" test row is dropped and the role
# progression shifts up by one row, with a new "Administrator" role appended.
$wsModifyUser = $wb.Worksheets.Item("modifyUser")
$wsModifyUser.Range("A2").Value = "modifyUser_Role"
$wsModifyUser.Range("C2").Value = "Mitarbeiter"
$wsModifyUser.Range("C3").Value = "Abonnent"
$wsModifyUser.Range("C4").Value = "Autor"
$wsModifyUser.Range("C5").Value = "Redakteur"
$wsModifyUser.Range("C6").Value = "Administrator"

# User sheet: row 6's role is now "Administrator" instead of "Mitarbeiter".
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Range("H6").Value = "Administrator"

# ---------------------------------------------------------------------------
# 2. Remove the deleteAdmin sheet entirely.
# ---------------------------------------------------------------------------
$excel.DisplayAlerts = $false
$wsDeleteAdmin = $wb.Worksheets.Item("deleteAdmin")
$wsDeleteAdmin.Delete()

# ---------------------------------------------------------------------------
# 3. Column width tweak on HTTP_LINK.
# ---------------------------------------------------------------------------
$wsHttpLink = $wb.Worksheets.Item("HTTP_LINK")
$wsHttpLink.Columns.Item(1).ColumnWidth = 41.3

# ---------------------------------------------------------------------------
# 4. Selection / active-sheet bookkeeping. Order matters: the last sheet we
#    select on becomes the workbook's active tab, so "User" goes last to end
#    up as activeTab.
# ---------------------------------------------------------------------------
$wsUserLogin = $wb.Worksheets.Item("UserLogin")
$wsUserLogin.Range("F2:F4").Select()

$wsModifyUser.Range("C6").Select()

$wsUser.Range("E12").Select()
